$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 data: product becomes "Cement" ---
$ws.Range("B2").Value = "Cement"
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 2.0
$ws.Range("E2").Value = 130.0
$ws.Range("F2").Value = 49.92
$ws.Range("G2").Value = 208.0
$ws.Range("H2").Value = 374.4

# --- Update row 3 data: product becomes "Supremo Beer 1L" ---
$ws.Range("B3").Value = "Supremo Beer 1L"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 1.0
$ws.Range("E3").Value = 150.0
$ws.Range("F3").Value = 19.8
$ws.Range("G3").Value = 165.0
$ws.Range("H3").Value = 156.75

# --- Remove the now-unused rows 4 and 5 ---
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# --- Recompute "best fit" column widths for the columns whose content changed ---
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(5).AutoFit()

# --- Refresh selection to the new used range ---
$ws.Range("A1:H3").Select()

Write-Host "Edit applied"
